# "We want data tables from dev" -- append the new httk 2.5.0 benchmark
# row to Sheet1, grow Table1 to cover it, and leave the selection where
# the user would have ended up after typing the last cell of the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 27

$ws.Cells.Item($newRow, 1).Value  = "2.5.0"
$ws.Cells.Item($newRow, 2).Value  = 1021
$ws.Cells.Item($newRow, 3).Value  = 1
$ws.Cells.Item($newRow, 4).Value  = 1
$ws.Cells.Item($newRow, 5).Value  = 0.9999
$ws.Cells.Item($newRow, 6).Value  = 0.9477
$ws.Cells.Item($newRow, 7).Value  = 353
$ws.Cells.Item($newRow, 8).Value  = 0.2716
$ws.Cells.Item($newRow, 9).Value  = 353
$ws.Cells.Item($newRow, 10).Value = 1.508
$ws.Cells.Item($newRow, 11).Value = 36
$ws.Cells.Item($newRow, 12).Value = 0.9698
$ws.Cells.Item($newRow, 13).Value = 80
$ws.Cells.Item($newRow, 14).Value = 1.132
$ws.Cells.Item($newRow, 15).Value = 80
$ws.Cells.Item($newRow, 16).Value = 0.6466
$ws.Cells.Item($newRow, 17).Value = 863
$ws.Cells.Item($newRow, 18).Value = "Added models 3comp2 and sumclearances"

# Match the left-aligned style ("s=1") already used by the rest of the table.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item($newRow, $col).HorizontalAlignment = -4131
}

# Grow the structured table to include the new row.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:R" + $newRow))

# Leave the selection on the last cell the user would have typed into.
$ws.Range("R" + $newRow).Select()
